$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-9 down to 3-10)
$ws.Rows(2).Insert()

# Populate the new row 2 with the "cruise" attribute metadata
$ws.Range("A2").Value = "cruise"
$ws.Range("B2").Value = "Identifier for research cruise generally including abbreviation for research vessel and voyage number"
$ws.Range("C2").Value = "character"

# The inserted row copies formatting from the row above (header row 1, which
# has wrapped-text style on column B); the target state has no explicit
# style on the new row, so clear the copied format from that cell only.
$ws.Range("B2").ClearFormats()

# Update the view/selection to match the committed state: whole row 2 selected
$ws.Range("A2:XFD2").Select()
